$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 124.9407424926758
$ws.Range("B3").Value = 131.4092559814453
$ws.Range("B4").Value = 128.3921051025391
$ws.Range("B5").Value = 123.5036163330078
$ws.Range("B6").Value = 125.4058151245117
$ws.Range("B7").Value = 123.9601135253906
$ws.Range("B8").Value = 126.0068969726562
$ws.Range("B9").Value = 122.9411087036133
$ws.Range("B10").Value = 126.1399612426758
$ws.Range("B11").Value = 125.5306701660156
$ws.Range("B12").Value = 128.0505065917969
$ws.Range("B13").Value = 136.166015625
$ws.Range("B14").Value = 140.05126953125
$ws.Range("B15").Value = 148.3330993652344
$ws.Range("B16").Value = 169.8754425048828
$ws.Range("B17").Value = 200.8213958740234
$ws.Range("B18").Value = 189.2234954833984
$ws.Range("B19").Value = 201.5697021484375
$ws.Range("B20").Value = 200.8118438720703
$ws.Range("B21").Value = 202.7599792480469
$ws.Range("B22").Value = 203.4048156738281
$ws.Range("B23").Value = 196.9842071533203
$ws.Range("B24").Value = 198.1509552001953
$ws.Range("B25").Value = 197.3408966064453
$ws.Range("B26").Value = 192.8445892333984
$ws.Range("B27").Value = 195.2014007568359
$ws.Range("B28").Value = 192.3062438964844
$ws.Range("B29").Value = 186.1914672851562
$ws.Range("B30").Value = 189.1931762695312
$ws.Range("B31").Value = 193.2124481201172
$ws.Range("B32").Value = 198.7594757080078
$ws.Range("B33").Value = 224.3966522216797
$ws.Range("B34").Value = 209.4432067871094
$ws.Range("B35").Value = 246.4754486083984
$ws.Range("B36").Value = 242.4913024902344
$ws.Range("B37").Value = 246.7333221435547
$ws.Range("B38").Value = 223.8305206298828
$ws.Range("B39").Value = 209.7524566650391
$ws.Range("B40").Value = 207.6990051269531
$ws.Range("B41").Value = 190.8677215576172
$ws.Range("B42").Value = 181.4691772460938
$ws.Range("B43").Value = 160.7430114746094
$ws.Range("B44").Value = 167.5046234130859
$ws.Range("B45").Value = 154.0199279785156
$ws.Range("B46").Value = 159.3650207519531
$ws.Range("B47").Value = 147.1376342773438
$ws.Range("B48").Value = 155.3755187988281
$ws.Range("B49").Value = 154.7249298095703
